$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.284.68"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").Value = "1.871.70"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4400"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3695"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9390"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "1.882.88"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.736"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.457"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06865"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "82.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("E18").Value = "  -4.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").Value = "28.270.29"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.139"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "2.123.59"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.027"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.340"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.733"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09050"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8011"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.852"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.44%  "
$ws.Range("E34").Value = "  -5.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.908"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.125"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01957"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.921"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.127"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5264"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1682"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.779"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06758"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4887"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.51%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000002455"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.684"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.07%  "
